$d = $word.ActiveDocument

# 1) Remove the stray _GoBack bookmark on the "Site :" paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2) "Réseaux sociaux ok ?" -> "Pas de réseaux sociaux -> Mettre dans cahier des charges"
$d.Content.Find.Execute("Réseaux sociaux ok ?", $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         "Pas de réseaux sociaux -> Mettre dans cahier des charges", 2)

# 3) Insert two new bullet items after the paragraph we just edited, reusing its
#    list formatting (ListParagraph style / numId 1 / ind left 360).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Pas de réseaux sociaux*") {
        $target = $p
    }
}

$target.Range.InsertParagraphAfter()
$photosPara = $target.Next()
$photosPara.Range.Text = "Sécuriser les photos"

$photosPara.Range.InsertParagraphAfter()
$loginPara = $photosPara.Next()
# Use a unique sentinel suffix so we can anchor a Find range precisely at the end
# of the new run (zero-length Bookmarks.Add only resolves cleanly against a real
# run boundary produced by Find), then add the bookmark there and erase the
# sentinel, leaving an empty _GoBack bookmark right after the text.
$loginPara.Range.Text = "Sécuriser les champs des formulaires + NO BLANK pour login admin@@MARK@@"

$sentinel = $d.Content
$sentinel.Find.Execute("@@MARK@@", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$d.Bookmarks.Add("_GoBack", $sentinel)
$sentinel.Text = ""
